$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Update header labels: "bandwidth(words/s)" -> "bandwidth(GB/s)" ---
$ws1.Range("P2").Value2 = "bandwidth(GB/s)"
$ws1.Range("Q2").Value2 = "bandwidth(GB/s)"
$ws1.Range("R2").Value2 = "bandwidth(GB/s)"
$ws1.Range("S2").Value2 = "bandwidth(GB/s)"

$ws2.Range("I2").Value2 = "bandwidth(GB/s)"
$ws2.Range("J2").Value2 = "bandwidth(GB/s)"

# --- Update bandwidth formulas on sheet 1 (Лист5), rows 3-51 ---
# old: =$T{r}/{col}{r}   (words/s, dividing memory-access count by time in ms)
# new: =$T{r}*4/({col}{r}/1000)/10^9  (GB/s, *4 bytes/word, ms->s, bytes->GB)
for ($r = 3; $r -le 51; $r++) {
    $ws1.Range("P$r").Formula = "=`$T$r*4/(D$r/1000)/10^9"
    $ws1.Range("Q$r").Formula = "=`$T$r*4/(G$r/1000)/10^9"
    $ws1.Range("R$r").Formula = "=`$T$r*4/(J$r/1000)/10^9"
    $ws1.Range("S$r").Formula = "=`$T$r*4/(M$r/1000)/10^9"
}

# --- Update bandwidth formulas on sheet 2 (Лист6), rows 3-67 ---
for ($r = 3; $r -le 67; $r++) {
    $ws2.Range("I$r").Formula = "=K$r*4/(D$r/1000)/10^9"
    $ws2.Range("J$r").Formula = "=K$r*4/(F$r/1000)/10^9"
}

# --- Update chart titles: "bandwidth" -> "bandwidth GB/s" ---
$chart2 = $ws1.ChartObjects().Item(2).Chart
$chart2.ChartTitle.Text = "Брзина (bandwidth GB/s)"

$chart4 = $ws2.ChartObjects().Item(2).Chart
$chart4.ChartTitle.Text = "Брзина (bandwidth GB/s)"

# --- Restore view / selection state to match the edited workbook ---
$ws1.Range("P2:S2").Select()
$ws2.Range("J3").Select()
